$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "44.125.84"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.366.41"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.21%  "
Set-TextValue "D5" "0.688"
$ws.Range("E5").Value = "  +2.07%  "
Set-TextValue "D6" "240.82"
$ws.Range("E6").Value = "  +2.19%  "
Set-TextValue "D7" "75.94"
$ws.Range("E7").Value = "  +4.47%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue "D9" "0.619"
$ws.Range("E9").Value = "  +8.74%  "
$ws.Range("E10").Value = "  +3.68%  "
Set-TextValue "D11" "57.09"
$ws.Range("E11").Value = "  +0.10%  "
Set-TextValue "D12" "32.98"
$ws.Range("E12").Value = "  +18.32%  "
Set-TextValue "D13" "7.40"
$ws.Range("E13").Value = "  +9.30%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "2.717.41"
$ws.Range("E15").Value = "  +0.62%  "
Set-TextValue "D16" "16.70"
$ws.Range("E16").Value = "  -1.01%  "
Set-TextValue "D17" "0.919"
$ws.Range("E17").Value = "  +3.72%  "
$ws.Range("D18").Value = "2.366.87"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "44.009.86"
$ws.Range("E19").Value = "  +0.30%  "
Set-TextValue "D20" "0.0000103"
$ws.Range("E20").Value = "  +2.37%  "
Set-TextValue "D21" "6.71"
$ws.Range("E21").Value = "  +5.70%  "
Set-TextValue "D22" "77.62"
$ws.Range("E22").Value = "  +1.16%  "
Set-TextValue "D23" "258.93"
$ws.Range("E23").Value = "  +2.18%  "
Set-TextValue "D25" "3.73"
$ws.Range("E25").Value = "  -1.52%  "
Set-TextValue "D26" "2.53"
$ws.Range("E26").Value = "  +1.54%  "
Set-TextValue "D27" "1.82"
$ws.Range("E27").Value = "  +15.38%  "
Set-TextValue "D28" "10.85"
$ws.Range("E28").Value = "  +3.37%  "
Set-TextValue "D29" "23.09"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("E30").Value = "  -1.55%  "
Set-TextValue "D31" "174.75"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  -1.55%  "
Set-TextValue "D33" "0.139"
$ws.Range("E33").Value = "  +4.34%  "
Set-TextValue "D34" "0.0766"
$ws.Range("E34").Value = "  +6.63%  "
Set-TextValue "D35" "5.38"
$ws.Range("E35").Value = "  +3.48%  "
Set-TextValue "D36" "5.43"
$ws.Range("E36").Value = "  +4.04%  "
Set-TextValue "D37" "3.77"
$ws.Range("E37").Value = "  +0.26%  "
Set-TextValue "D38" "2.39"
$ws.Range("E38").Value = "  -1.71%  "
Set-TextValue "D39" "6.40"
Set-TextValue "D40" "0.0285"
$ws.Range("E40").Value = "  +3.66%  "
Set-TextValue "D41" "0.211"
$ws.Range("E41").Value = "  +14.32%  "
Set-TextValue "D42" "19.77"
$ws.Range("E42").Value = "  +3.27%  "
Set-TextValue "D43" "0.110"
$ws.Range("E43").Value = "  +12.45%  "
Set-TextValue "D44" "9.27"
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("E45").Value = "  -0.02%  "
Set-TextValue "D46" "4.85"
$ws.Range("E46").Value = "  +9.18%  "
$ws.Range("E47").Value = "  +10.66%  "
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("E49").Value = "  +2.31%  "
Set-TextValue "D50" "101.48"
$ws.Range("E50").Value = "  +4.09%  "
Set-TextValue "D51" "56.62"
$ws.Range("E51").Value = "  +8.47%  "
